# Stage at DAM corrected
# - Adds a new "DAM" worksheet (after "LBJ") with stage-correction data
# - Updates LBJ's column A width/selection to reflect the new active tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the existing "LBJ" sheet: widen column A, clear its old selection
# ---------------------------------------------------------------------------
$lbj = $wb.Worksheets.Item("LBJ")

# Column A: 12.71 (bestFit) -> 22 (explicit custom width, no bestFit)
$lbj.Columns.Item(1).ColumnWidth = 21.14

# Select the whole of column A (this also clears the old F9 selection)
$lbj.Columns.Item(1).Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Add the new "DAM" worksheet right after "LBJ"
# ---------------------------------------------------------------------------
$dam = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lbj)
$dam.Name = "DAM"

# Headers (row 1) - reuse the same shared strings as LBJ
$dam.Range("A1").Value = "T1_date"
$dam.Range("B1").Value = "T1_time"
$dam.Range("C1").Value = "T2_date"
$dam.Range("D1").Value = "T2_time"
$dam.Range("E1").Value = "z"

# Data rows (2-6)
$dam.Range("A2").Value = 41394.59375
$dam.Range("B2").Value = 1415
$dam.Range("C2").Value = 41430.5625
$dam.Range("D2").Value = 1330
$dam.Range("E2").Value = 3

$dam.Range("A3").Value = 41629.53125
$dam.Range("B3").Value = 1245
$dam.Range("C3").Value = 41639.489583333336
$dam.Range("D3").Value = 1145
$dam.Range("E3").Value = 1

$dam.Range("A4").Value = 41727.666666666664
$dam.Range("B4").Value = 1600
$dam.Range("C4").Value = 41807.145833333336
$dam.Range("D4").Value = 330
$dam.Range("E4").Value = 0.5

$dam.Range("A5").Value = 41807.15625
$dam.Range("B5").Value = 345
$dam.Range("C5").Value = 41899.208333333336
$dam.Range("D5").Value = 500
$dam.Range("E5").Value = 1

$dam.Range("A6").Value = 41899.21875
$dam.Range("B6").Value = 515
$dam.Range("C6").Value = 41903.635416666664
$dam.Range("D6").Value = 1515
$dam.Range("E6").Value = 1.5

# Column widths matching the LBJ sheet layout
$dam.Columns.Item(1).ColumnWidth = 21.14
$dam.Columns.Item(2).ColumnWidth = $lbj.Columns.Item(2).ColumnWidth
$dam.Columns.Item(3).ColumnWidth = $lbj.Columns.Item(3).ColumnWidth
$dam.Columns.Item(4).ColumnWidth = $lbj.Columns.Item(4).ColumnWidth

# Number formats: columns A/C are dates, B/D are 4-digit times
$dam.Columns.Item(1).NumberFormat = "m/d/yy\ h:mm;@"
$dam.Columns.Item(2).NumberFormat = "0000"
$dam.Columns.Item(3).NumberFormat = "m/d/yy\ h:mm;@"
$dam.Columns.Item(4).NumberFormat = "0000"

# Page setup (portrait, matching LBJ)
$dam.PageSetup.Orientation = 1

# Selection/active cell on the new sheet
$dam.Range("A7").Select() | Out-Null
